# ---------------------------------------------------------------------------
# OLX Monitor 2026-02-24 09:38
# Append the 8 newly-discovered listings (rows 219-226) to the running log
# kept on the "PODSUMOWANIE" worksheet, mirroring the formatting used by the
# existing log rows (row 218 is the most recent pre-existing entry).
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PODSUMOWANIE")

# 1) Seed the new rows with the same cell formatting as the last existing row.
$ws.Range("A218:H218").Copy($ws.Range("A219:H226"))

# 2) Rows 220 and 223 are "fresh" listings that use the alternate (non-highlighted)
#    style for the "days listed" column; pull that style from row 7, which already
#    uses it, without touching its value.
$ws.Range("F7").Copy($ws.Range("F220"))
$ws.Range("F7").Copy($ws.Range("F223"))

# 3) Fill in the actual data for each new row.

# Row 219: poqui - Mieszkanie z KLIMATYZACJĄ 5 minut od UMCS, UP, KUL - Długosza
$ws.Cells.Item(219, 1).Value = "2026-02-24 09:38:34"
$ws.Cells.Item(219, 2).Value = "poqui"
$ws.Cells.Item(219, 3).Value = "Mieszkanie z KLIMATYZACJĄ 5 minut od UMCS, UP, KUL - Długosza"
$ws.Cells.Item(219, 4).Value = 2049
$ws.Cells.Item(219, 5).Value = "19.12.2025"
$ws.Cells.Item(219, 6).Value = 66
$ws.Cells.Item(219, 7).Value = "https://www.olx.pl/d/oferta/mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc.html"
$ws.Cells.Item(219, 8).Value = "mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc"

# Row 220: poqui - Świeżo wykończone mieszkanie z dużym balkonem - Ponikwoda
$ws.Cells.Item(220, 1).Value = "2026-02-24 09:38:34"
$ws.Cells.Item(220, 2).Value = "poqui"
$ws.Cells.Item(220, 3).Value = "Świeżo wykończone mieszkanie z dużym balkonem - Ponikwoda"
$ws.Cells.Item(220, 4).Value = 2299
$ws.Cells.Item(220, 5).Value = "19.01.2026"
$ws.Cells.Item(220, 6).Value = 35
$ws.Cells.Item(220, 7).Value = "https://www.olx.pl/d/oferta/swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR.html"
$ws.Cells.Item(220, 8).Value = "swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR"

# Row 221: poqui - Kawalerka po remoncie z funkcjonalną antresolą - ul. Jana Sawy
$ws.Cells.Item(221, 1).Value = "2026-02-24 09:38:34"
$ws.Cells.Item(221, 2).Value = "poqui"
$ws.Cells.Item(221, 3).Value = "Kawalerka po remoncie z funkcjonalną antresolą - ul. Jana Sawy"
$ws.Cells.Item(221, 4).Value = 2499
$ws.Cells.Item(221, 5).Value = "28.10.2025"
$ws.Cells.Item(221, 6).Value = 118
$ws.Cells.Item(221, 7).Value = "https://www.olx.pl/d/oferta/kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger.html"
$ws.Cells.Item(221, 8).Value = "kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger"

# Row 222: poqui - Przytulny pokój blisko Politechniki – ul. Przytulna
$ws.Cells.Item(222, 1).Value = "2026-02-24 09:38:34"
$ws.Cells.Item(222, 2).Value = "poqui"
$ws.Cells.Item(222, 3).Value = "Przytulny pokój blisko Politechniki – ul. Przytulna"
$ws.Cells.Item(222, 4).Value = 549
$ws.Cells.Item(222, 5).Value = "'10.10.2025"
$ws.Cells.Item(222, 6).Value = 137
$ws.Cells.Item(222, 7).Value = "https://www.olx.pl/d/oferta/przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz.html"
$ws.Cells.Item(222, 8).Value = "przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz"

# Row 223: pokojewlublinie - WOLNY OD ZARAZ! Super lokalizacja, blisko centrum, ul. Paganiniego 12
$ws.Cells.Item(223, 1).Value = "2026-02-24 09:38:34"
$ws.Cells.Item(223, 2).Value = "pokojewlublinie"
$ws.Cells.Item(223, 3).Value = "WOLNY OD ZARAZ! Super lokalizacja, blisko centrum, ul. Paganiniego 12"
$ws.Cells.Item(223, 4).Value = 12640
$ws.Cells.Item(223, 5).Value = "19.01.2026"
$ws.Cells.Item(223, 6).Value = 35
$ws.Cells.Item(223, 7).Value = "https://www.olx.pl/d/oferta/wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc.html"
$ws.Cells.Item(223, 8).Value = "wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc"

# Row 224: pokojewlublinie - WOLNY OD ZARAZ! Pokój jedynka, ul. Romanowskiego 58
$ws.Cells.Item(224, 1).Value = "2026-02-24 09:38:34"
$ws.Cells.Item(224, 2).Value = "pokojewlublinie"
$ws.Cells.Item(224, 3).Value = "WOLNY OD ZARAZ! Pokój jedynka, ul. Romanowskiego 58"
$ws.Cells.Item(224, 4).Value = 0
$ws.Cells.Item(224, 5).Value = "'11.08.2025"
$ws.Cells.Item(224, 6).Value = 196
$ws.Cells.Item(224, 7).Value = "https://www.olx.pl/d/oferta/wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm.html"
$ws.Cells.Item(224, 8).Value = "wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm"

# Row 225: dawnypatron - Ładny pokój jednoosobowy. Wynajmę duży pokój w centrum. ul Niecała 4.
$ws.Cells.Item(225, 1).Value = "2026-02-24 09:38:34"
$ws.Cells.Item(225, 2).Value = "dawnypatron"
$ws.Cells.Item(225, 3).Value = "Ładny pokój jednoosobowy. Wynajmę duży pokój w centrum. ul Niecała 4."
$ws.Cells.Item(225, 4).Value = 730
$ws.Cells.Item(225, 5).Value = "20.09.2024"
$ws.Cells.Item(225, 6).Value = 521
$ws.Cells.Item(225, 7).Value = "https://www.olx.pl/d/oferta/ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM.html"
$ws.Cells.Item(225, 8).Value = "ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM"

# Row 226: dawnypatron - Mam do wynajęcia pokój dla os. pracującej lub studenta. Narutowicza 14
$ws.Cells.Item(226, 1).Value = "2026-02-24 09:38:34"
$ws.Cells.Item(226, 2).Value = "dawnypatron"
$ws.Cells.Item(226, 3).Value = "Mam do wynajęcia pokój dla os. pracującej lub studenta. Narutowicza 14"
$ws.Cells.Item(226, 4).Value = 14690
$ws.Cells.Item(226, 5).Value = "'05.12.2025"
$ws.Cells.Item(226, 6).Value = 80
$ws.Cells.Item(226, 7).Value = "https://www.olx.pl/d/oferta/mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv.html"
$ws.Cells.Item(226, 8).Value = "mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv"

# 4) The apostrophe trick above flags those cells with a quote-prefix style;
#    restore the plain "date string" formatting (style of column E) used by the
#    rest of the log, without touching the text we just entered.
$ws.Range("E7").Copy()
$ws.Range("E222").PasteSpecial(-4122)
$ws.Range("E224").PasteSpecial(-4122)
$ws.Range("E226").PasteSpecial(-4122)
$excel.CutCopyMode = 0

